$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 44 and 45 swapped coins (TheGraph <-> Stacks), including their updated
# price / volume figures.
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'2.93"
$ws.Range("E44").Value = "  -8.33%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.319"
$ws.Range("E45").Value = "  -6.48%  "

# Updated Price (column D) and Volume(1h) (column E) figures for the remaining
# rows. A leading apostrophe is used for values that would otherwise be parsed
# by Excel as numbers, so that they are stored as text like the rest of the
# sheet (mirroring the original "Price"/"Volume" columns, which are textual).
$ws.Range("D2").Value = "66.513.89"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.247.71"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'578.27"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "'170.29"
$ws.Range("E6").Value = "  -8.33%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").Value = "3.242.86"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("E10").Value = "  -6.79%  "
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").Value = "'44.31"
$ws.Range("E12").Value = "  -5.73%  "
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "'682.51"
$ws.Range("E14").Value = "  +4.23%  "
$ws.Range("D15").Value = "3.780.59"
$ws.Range("E15").Value = "  +4.21%  "
$ws.Range("D16").Value = "'8.14"
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("D17").Value = "66.550.99"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "3.253.80"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").Value = "'17.00"
$ws.Range("E20").Value = "  -4.97%  "
$ws.Range("D21").Value = "'10.52"
$ws.Range("E21").Value = "  -5.53%  "
$ws.Range("D22").Value = "'0.871"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").Value = "'16.72"
$ws.Range("E23").Value = "  -5.55%  "
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").Value = "'96.64"
$ws.Range("E25").Value = "  -3.65%  "
$ws.Range("E26").Value = "  -4.82%  "
$ws.Range("D27").Value = "'2.60"
$ws.Range("E27").Value = "  -6.94%  "
$ws.Range("D28").Value = "'8.84"
$ws.Range("E28").Value = "  -8.08%  "
$ws.Range("D29").Value = "'32.14"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'8.18"
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("D31").Value = "'6.58"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("D32").Value = "'569.08"
$ws.Range("E32").Value = "  -5.35%  "
$ws.Range("D33").Value = "'10.73"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").Value = "3.771.47"
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -4.45%  "
$ws.Range("D37").Value = "'54.85"
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("E38").Value = "  -16.93%  "
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("E40").Value = "  -7.61%  "
$ws.Range("D41").Value = "'30.99"
$ws.Range("E41").Value = "  -6.61%  "
$ws.Range("D42").Value = "'3.23"
$ws.Range("E42").Value = "  -4.05%  "
$ws.Range("D43").Value = "0.0₃0644"
$ws.Range("E43").Value = "  -8.16%  "
$ws.Range("D46").Value = "'0.0396"
$ws.Range("E46").Value = "  -4.98%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D49").Value = "'2.49"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").Value = "'1.32"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "'126.35"
$ws.Range("E51").Value = "  -2.93%  "
